$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 9 (current row 9 "Residential NG" and everything
# below shifts down by one to make room for the new "Residential solar
# thermal" entry).
$ws.Rows.Item(9).Insert()

# Row 8 ("Residential kerosene") now points at the renamed OSEMOSYS code.
$ws.Range("B8").Value() = "RES_CWH_KER_001"

# Populate the newly inserted row 9 with the new dataset entry.
$ws.Range("A9").Value() = "Residential solar thermal"
$ws.Range("B9").Value() = "RES_CWH_SOLAR"
$ws.Range("C9").Value() = "Residential"
$ws.Range("D9").Value() = "Water heating mainly"
$ws.Range("E9").Value() = "Solar"

# Restore the selection to match the authored workbook.
$ws.Range("A8").Select()
